$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 6619.1816
$ws.Range("I74").Value = 9448.235000000001
$ws.Range("J74").Value = 3613.3125
$ws.Range("K74").Value = 9448.235000000001
$ws.Range("L74").Value = 3613.3125
$ws.Range("M74").Value = -8512.235000000001
$ws.Range("N74").Value = -5485.3125
$ws.Range("H76").Value = 8467.125
$ws.Range("I76").Value = 10677.117
$ws.Range("J76").Value = 3100
$ws.Range("K76").Value = 10677.117
$ws.Range("L76").Value = 3100
$ws.Range("M76").Value = -10362.117
$ws.Range("N76").Value = -3730
$ws.Range("H77").Value = 6619.1816
$ws.Range("I77").Value = 9448.235000000001
$ws.Range("J77").Value = 3613.3125
$ws.Range("K77").Value = 47241.175
$ws.Range("L77").Value = 18066.5625
$ws.Range("M77").Value = -42561.175
$ws.Range("N77").Value = -27426.5625
$ws.Range("H79").Value = 8467.125
$ws.Range("I79").Value = 10677.117
$ws.Range("J79").Value = 3100
$ws.Range("K79").Value = 10677.117
$ws.Range("L79").Value = 3100
$ws.Range("M79").Value = -9585.117
$ws.Range("N79").Value = -5284
$ws.Range("H101").Value = 1276.2941
$ws.Range("I101").Value = 1108.0834
$ws.Range("K101").Value = 3324.2502
$ws.Range("M101").Value = -1702.2502
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 80009
$ws.Range("J17").Value = 80009
$ws.Range("L17").Value = 80009
$ws.Range("N17").Value = -80355
$ws.Range("H45").Value = 3263.7827
$ws.Range("I45").Value = 2682.3125
$ws.Range("K45").Value = 2682.3125
$ws.Range("M45").Value = -2305.3125
$ws.Range("H80").Value = 41779.6
$ws.Range("J80").Value = 41779.6
$ws.Range("L80").Value = 41779.6
$ws.Range("N80").Value = -43775.6
$ws.Range("H83").Value = 41779.6
$ws.Range("J83").Value = 41779.6
$ws.Range("L83").Value = 125338.8
$ws.Range("N83").Value = -135322.8
$ws.Range("H117").Value = 22147.666
$ws.Range("J117").Value = 22147.666
$ws.Range("L117").Value = 22147.666
$ws.Range("N117").Value = -31325.666
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 30405
$ws.Range("I82").Value = 11525.667
$ws.Range("J82").Value = 34761.77
$ws.Range("K82").Value = 11525.667
$ws.Range("L82").Value = 34761.77
$ws.Range("M82").Value = -11142.667
$ws.Range("N82").Value = -35527.77
$ws.Range("H85").Value = 30405
$ws.Range("I85").Value = 11525.667
$ws.Range("J85").Value = 34761.77
$ws.Range("K85").Value = 11525.667
$ws.Range("L85").Value = 34761.77
$ws.Range("M85").Value = -10199.667
$ws.Range("N85").Value = -37413.77
$ws.Range("H86").Value = 1802.5
$ws.Range("I86").Value = 1907.125
$ws.Range("J86").Value = 1174.75
$ws.Range("K86").Value = 1907.125
$ws.Range("L86").Value = 1174.75
$ws.Range("M86").Value = -784.125
$ws.Range("N86").Value = -3420.75
$ws.Range("H89").Value = 1802.5
$ws.Range("I89").Value = 1907.125
$ws.Range("J89").Value = 1174.75
$ws.Range("K89").Value = 9535.625
$ws.Range("L89").Value = 5873.75
$ws.Range("M89").Value = -3919.625
$ws.Range("N89").Value = -17105.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 14926.667
$ws.Range("J41").Value = 16932
$ws.Range("L41").Value = 16932
$ws.Range("N41").Value = -17788
$ws.Range("H50").Value = 9156.6
$ws.Range("J50").Value = 9156.6
$ws.Range("L50").Value = 9156.6
$ws.Range("N50").Value = -10406.6
$ws.Range("H60").Value = 36925
$ws.Range("J60").Value = 36925
$ws.Range("L60").Value = 36925
$ws.Range("N60").Value = -37947
$ws.Range("H68").Value = 16802.75
$ws.Range("J68").Value = 16802.75
$ws.Range("L68").Value = 16802.75
$ws.Range("N68").Value = -18300.75
$ws.Range("H71").Value = 16802.75
$ws.Range("J71").Value = 16802.75
$ws.Range("L71").Value = 50408.25
$ws.Range("N71").Value = -57896.25
$ws.Range("H109").Value = 11000
$ws.Range("J109").Value = 11000
$ws.Range("L109").Value = 11000
$ws.Range("N109").Value = -13080
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 7360.1
$ws.Range("I80").Value = 4350.5
$ws.Range("J80").Value = 8112.5
$ws.Range("K80").Value = 13051.5
$ws.Range("L80").Value = 24337.5
$ws.Range("M80").Value = -12115.5
$ws.Range("N80").Value = -26209.5
$ws.Range("H83").Value = 7360.1
$ws.Range("I83").Value = 4350.5
$ws.Range("J83").Value = 8112.5
$ws.Range("K83").Value = 39154.5
$ws.Range("L83").Value = 73012.5
$ws.Range("M83").Value = -34474.5
$ws.Range("N83").Value = -82372.5
$ws.Range("H136").Value = 3445.7693
$ws.Range("I136").Value = 1343.3334
$ws.Range("K136").Value = 4030.0002
$ws.Range("M136").Value = 1069.9998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 17527.1
$ws.Range("J57").Value = 18912.889
$ws.Range("L57").Value = 18912.889
$ws.Range("N57").Value = -20552.889
$ws.Range("H80").Value = 3316
$ws.Range("I80").Value = 2770
$ws.Range("J80").Value = 5500
$ws.Range("K80").Value = 2770
$ws.Range("L80").Value = 5500
$ws.Range("M80").Value = -1772
$ws.Range("N80").Value = -7496
$ws.Range("H83").Value = 3316
$ws.Range("I83").Value = 2770
$ws.Range("J83").Value = 5500
$ws.Range("K83").Value = 13850
$ws.Range("L83").Value = 27500
$ws.Range("M83").Value = -8858
$ws.Range("N83").Value = -37484
$ws.Range("H107").Value = 612.0741
$ws.Range("I107").Value = 439.52942
$ws.Range("J107").Value = 905.4
$ws.Range("K107").Value = 439.52942
$ws.Range("L107").Value = 905.4
$ws.Range("M107").Value = 1480.47058
$ws.Range("N107").Value = -4745.4
$ws.Range("H123").Value = 34575.332
$ws.Range("J123").Value = 34575.332
$ws.Range("L123").Value = 34575.332
$ws.Range("N123").Value = -39475.332
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2858.261
$ws.Range("I122").Value = 2217.1428
$ws.Range("J122").Value = 3855.5557
$ws.Range("K122").Value = 6651.428400000001
$ws.Range("L122").Value = 11566.6671
$ws.Range("M122").Value = -4201.428400000001
$ws.Range("N122").Value = -16466.6671
$ws.Range("H138").Value = 46984.11
$ws.Range("J138").Value = 46984.11
$ws.Range("L138").Value = 46984.11
$ws.Range("N138").Value = -57264.11
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 18777
$ws.Range("J109").Value = 18777
$ws.Range("L109").Value = 18777
$ws.Range("N109").Value = -21551
$ws.Range("H135").Value = 62327.855
$ws.Range("J135").Value = 66382.5
$ws.Range("L135").Value = 66382.5
$ws.Range("N135").Value = -76522.5
